$wb = $excel.ActiveWorkbook

# --- Sheet: Folder Inventory ---
$ws = $wb.Worksheets.Item("Folder Inventory")

# Insert a new row at row 2 (shifts existing rows 2..68 down to 3..69)
$ws.Rows.Item(2).Insert()

# Insert() copies the formatting of the row above (the bold header row) -
# reset to the plain/default style used by the other data rows.
$ws.Rows.Item(2).Style = "Normal"

# Populate the newly inserted row with the new folder entry
$ws.Cells.Item(2, 1).Value = "Azure Local Hands-on Lab"
$ws.Cells.Item(2, 2).Value = "Azure Local Hands-on Lab"
$ws.Cells.Item(2, 3).Value = "2025-06-11 19:56:28 +0530"
$ws.Cells.Item(2, 4).Value = 1
$ws.Cells.Item(2, 5).Value = "Root"

# --- Sheet: Metadata ---
$meta = $wb.Worksheets.Item("Metadata")
$meta.Cells.Item(3, 2).Value = "2025-06-11 14:26:46 UTC"
$meta.Cells.Item(4, 2).Value = 68

# "Workflow Run" is stored as text ("2"), not a number - force text type,
# write the value, then restore the cell's original (unstyled) appearance.
$meta.Cells.Item(5, 2).NumberFormat = "@"
$meta.Cells.Item(5, 2).Value = "2"
$meta.Cells.Item(5, 2).Style = "Normal"

# --- Sheet: Summary ---
$summary = $wb.Worksheets.Item("Summary")
$summary.Cells.Item(2, 2).Value = 68
$summary.Cells.Item(3, 2).Value = 68
$summary.Cells.Item(5, 2).Value = "2025-06-11 19:56:28 +0530"
